# Updated symbol list on Tue Dec 27 19:33:49 UTC 2022 with GitHub Actions
#
# Refreshes the "Price" column with new quotes, bumps two "Volume(1h)"
# labels to reflect new Worst-in-24h / Best-in-24h markers, and re-syncs
# the BKEXToken / CEJI / KickToken block (rows 41-43) whose relative
# ranking shuffled in this run (KickToken moved up to 41, BKEXToken
# dropped to 42, CEJI dropped to 43), including their coin links.
#
# The Price column stores plain numeric-looking text (e.g. "0.05830",
# "245.72") rather than real numbers, so each Price cell is forced to
# Text format before the write -- otherwise Excel would silently coerce
# the digits to a Number and drop meaningful trailing zeros.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextPrice($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
}

Set-TextPrice "D2" "245.72"
Set-TextPrice "D4" "5.323"
Set-TextPrice "D5" "0.05830"
Set-TextPrice "D6" "6.473"
Set-TextPrice "D7" "3.348"
Set-TextPrice "D8" "0.8107"
Set-TextPrice "D9" "0.9222"
Set-TextPrice "D10" "0.1411"
Set-TextPrice "D11" "0.07354"
Set-TextPrice "D12" "0.03067"
Set-TextPrice "D13" "0.03074"
Set-TextPrice "D14" "0.09359"
Set-TextPrice "D15" "3.856"
Set-TextPrice "D16" "0.001559"
Set-TextPrice "D17" "0.04675"

Set-TextPrice "D18" "0.0005995"
$ws.Range("E18").Value = "17OneONEWorstin24h"

Set-TextPrice "D19" "0.006043"
Set-TextPrice "D20" "0.001242"
Set-TextPrice "D22" "0.00008806"
Set-TextPrice "D23" "3.592"
Set-TextPrice "D26" "0.1328"
Set-TextPrice "D28" "0.0002341"
Set-TextPrice "D40" "0.03840"

# Row 41: BKEXToken -> KickToken
$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
Set-TextPrice "D41" "0.006344"
$ws.Range("E41").Value = "40KickTokenKICK"

# Row 42: CEJI -> BKEXToken
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextPrice "D42" "0.1065"
$ws.Range("E42").Value = "41BKEXTokenBKK"

# Row 43: KickToken -> CEJI
$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextPrice "D43" "0.003202"
$ws.Range("E43").Value = "42CEJICEJI"

Set-TextPrice "D44" "0.007624"
Set-TextPrice "D45" "0.00005265"
Set-TextPrice "D46" "0.00000000750"
Set-TextPrice "D47" "0.6805"
Set-TextPrice "D48" "0.001835"
Set-TextPrice "D49" "0.00002101"
Set-TextPrice "D50" "0.0002001"
